$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add header cells for the new columns AD, AE, AF
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Copy the style of an existing header cell (AC1) so the new headers look consistent (bold, bordered, centered)
$ws.Range("AC1").Copy() | Out-Null
$ws.Range("AD1:AF1").PasteSpecial(-4122) | Out-Null  # xlPasteFormats

# Fill in team record values for each data row (2 through 45)
for ($r = 2; $r -le 45; $r++) {
    $ws.Cells.Item($r, 30).Value = 43   # AD = column 30
    $ws.Cells.Item($r, 31).Value = 119  # AE = column 31
    $ws.Cells.Item($r, 32).Value = 0    # AF = column 32
}

$excel.CutCopyMode = 0
